# Apply value updates to existing rows (2-6)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 1.1
$ws.Range("S2").Value = 1.05
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04

$ws.Range("AA3").Value = 1000

$ws.Range("G4").Value = 2.46
$ws.Range("I4").Value = 3.9
$ws.Range("J4").Value = 3.1
$ws.Range("T4").Value = 2.04
$ws.Range("V4").Value = 1.34
$ws.Range("X4").Value = 9.199999999999999
$ws.Range("AB4").Value = 8
$ws.Range("AD4").Value = 20
$ws.Range("AG4").Value = 15
$ws.Range("AH4").Value = 24
$ws.Range("AI4").Value = 90

$ws.Range("J5").Value = 6.8
$ws.Range("L5").Value = 1.14
$ws.Range("Q5").Value = 1.34
$ws.Range("R5").Value = 1.96

$ws.Range("H6").Value = 4.8
$ws.Range("O6").Value = 1.24
$ws.Range("P6").Value = 2.12
$ws.Range("T6").Value = 1.73
$ws.Range("U6").Value = 2.1
$ws.Range("AN6").Value = 9.6

# Append two new match rows (7 and 8) -- columns A..AO (41 values each)
$row7 = @("Argentinian Primera Division", "2025-11-14", "20:00:00", "Lanus", "Atl Tucuman", 1.84, 2.04, 4.6, 5.4, 3.1, 3.8, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$row8 = @("Brazilian Serie B", "2025-11-14", "20:00:00", "Paysandu", "Amazonas FC", 3.25, 3.7, 2.32, 2.52, 3.15, 3.6, 1.48, 1.09, 3.05, 1.42, 1.7, 2.2, 1.26, 4.2, 1.89, 1.93, 1.65, 1.37, 12, 9.199999999999999, 15.5, 980, 12, 7.8, 12.5, 32, 25, 16, 21, 980, 75, 50, 65, 160, 60, 29)

# Column B holds dates formatted as plain text like "2025-11-14"; force text
# formatting so Excel does not convert the string into a date serial number.
$ws.Range("B7:B8").NumberFormat = "@"

for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $row7[$i]
}
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $row8[$i]
}

# Restore the default "Normal" style on the date cells so no extra,
# text-specific style is persisted (matches the rest of the sheet which
# carries no explicit cell style on data rows).
$ws.Range("B7:B8").Style = "Normal"
